{"js": "// Apply Dutch translations to HIVE TEAMS.docx (Office.js / Word JavaScript API)\n\nconst body = context.document.body;\n\n// 1) Intro paragraph: replace the trailing English sentence with its Dutch translation.\nconst introResults = body.search(\n  \"After SmartCash starts to require more teams they will be created and these will be splintered into smaller teams.\",\n  { matchCase: true }\n);\nintroResults.load(\"text\");\nawait context.sync();\nif (introResults.items.length > 0) {\n  introResults.items[0].insertText(\n    \"Wanneer SmartCash de behoefte ziet ontstaan, zullen meer teams worden aangemaakt die weer verdeeld zullen worden in kleinere teams.\",\n    \"Replace\"\n  );\n}\n\n// 2) Remove the stray leading (non-breaking-space-only) run that sits right before\n//    the \"HIVE TEAM: OUTREACH\" heading text (not the \"OUTREACH 2\" heading).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\\u00a0HIVE TEAM: OUTREACH\") {\n    const nbspResults = paragraphs.items[i].search(\"\\u00a0\", { matchCase: true });\n    nbspResults.load(\"text\");\n    await context.sync();\n    if (nbspResults.items.length > 0) {\n      nbspResults.items[0].delete();\n    }\n    break;\n  }\n}\nawait context.sync();\n\n// 3) \"This team focuses on community building, growth, general user acquisition.\"\n//    (exact match only \u2014 must not touch the \"... in South America\" variant elsewhere.)\nconst teamFocusResults = body.search(\n  \"This team focuses on community building, growth, general user acquisition.\",\n  { matchCase: true }\n);\nteamFocusResults.load(\"text\");\nawait context.sync();\nif (teamFocusResults.items.length > 0) {\n  teamFocusResults.items[0].insertText(\n    \"Het team richt zich op het versterken van de community, aanjagen van groei en toevoegen van nieuwe gebruikers.\",\n    \"Replace\"\n  );\n}\n\n// 4) \"Best known for creating the original Dash Force proposal.\"\nconst dashForceResults = body.search(\n  \"Best known for creating the original Dash Force proposal.\",\n  { matchCase: true }\n);\ndashForceResults.load(\"text\");\nawait context.sync();\nif (dashForceResults.items.length > 0) {\n  dashForceResults.items[0].insertText(\n    \"Meest bekend voor het maken van het originele Dash Force-voorstel.\",\n    \"Replace\"\n  );\n}\n\n// 5) \"SmartCash Advisor\" -> \"SmartCash Adviseur\"\nconst advisorResults = body.search(\"SmartCash Advisor\", { matchCase: true });\nadvisorResults.load(\"text\");\nawait context.sync();\nif (advisorResults.items.length > 0) {\n  advisorResults.items[0].insertText(\"SmartCash Adviseur\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Apply Dutch translations to HIVE TEAMS.docx (Word COM interop)\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphTextByExactMatch($doc, [string]$oldText, [string]$newText) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $t = $p.Range.Text.TrimEnd([char]0x0D)\n        if ($t -eq $oldText) {\n            $p.Range.Text = $newText\n            return $true\n        }\n    }\n    return $false\n}\n\n# 1) Intro paragraph: replace the trailing English sentence with its Dutch translation\n#    (whole-paragraph replace keeps the run's xml:space=\"preserve\" intact).\nSet-ParagraphTextByExactMatch $d `\n    \"SmartCash streeft naar een gedecentraliseerd teamstructuur door op een effici\u00ebnte manier de werkdruk over meerdere wereldwijde Hive teams te verdelen. De teams hieronder zorgen er gewoon voor dat dingen op gang komen. After SmartCash starts to require more teams they will be created and these will be splintered into smaller teams.\" `\n    \"SmartCash streeft naar een gedecentraliseerd teamstructuur door op een effici\u00ebnte manier de werkdruk over meerdere wereldwijde Hive teams te verdelen. De teams hieronder zorgen er gewoon voor dat dingen op gang komen. Wanneer SmartCash de behoefte ziet ontstaan, zullen meer teams worden aangemaakt die weer verdeeld zullen worden in kleinere teams.\" `\n    | Out-Null\n\n# 2) Remove the stray leading (non-breaking-space-only) run that sits right before\n#    the \"HIVE TEAM: OUTREACH\" heading text (not the \"OUTREACH 2\" heading).\n$targetText = [char]0x00A0 + \"HIVE TEAM: OUTREACH\"\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $paraText = $p.Range.Text.TrimEnd([char]0x0D)\n    if ($paraText -eq $targetText) {\n        $r = $p.Range\n        $nbspRange = $d.Range($r.Start, $r.Start + 1)\n        $nbspRange.Delete()\n        break\n    }\n}\n\n# 3) \"This team focuses on community building, growth, general user acquisition.\"\n#    (exact match only \u2014 must not touch the \"... in South America\" variant elsewhere.)\nSet-ParagraphTextByExactMatch $d `\n    \"This team focuses on community building, growth, general user acquisition.\" `\n    \"Het team richt zich op het versterken van de community, aanjagen van groei en toevoegen van nieuwe gebruikers.\" `\n    | Out-Null\n\n# 4) \"Best known for creating the original Dash Force proposal.\"\nSet-ParagraphTextByExactMatch $d `\n    \"Best known for creating the original Dash Force proposal.\" `\n    \"Meest bekend voor het maken van het originele Dash Force-voorstel.\" `\n    | Out-Null\n\n# 5) \"SmartCash Advisor\" -> \"SmartCash Adviseur\"\nSet-ParagraphTextByExactMatch $d \"SmartCash Advisor\" \"SmartCash Adviseur\" | Out-Null\n"}
